$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '51.712.11'
Set-TextValue $ws.Range("E2") '  -1.21%  '
Set-TextValue $ws.Range("D3") '2.784.92'
Set-TextValue $ws.Range("E3") '  -1.94%  '
Set-TextValue $ws.Range("E4") '  +0.03%  '
Set-TextValue $ws.Range("D5") '358.32'
Set-TextValue $ws.Range("E5") '  -0.86%  '
Set-TextValue $ws.Range("D6") '109.78'
Set-TextValue $ws.Range("E6") '  -2.84%  '
Set-TextValue $ws.Range("D7") '0.554'
Set-TextValue $ws.Range("E7") '  -3.19%  '
Set-TextValue $ws.Range("E8") '  +0.03%  '
Set-TextValue $ws.Range("D9") '0.589'
Set-TextValue $ws.Range("E9") '  -2.81%  '
Set-TextValue $ws.Range("D10") '39.78'
Set-TextValue $ws.Range("E10") '  -3.54%  '
Set-TextValue $ws.Range("D11") '0.138'
Set-TextValue $ws.Range("E11") '  +4.31%  '
Set-TextValue $ws.Range("D12") '0.0844'
Set-TextValue $ws.Range("E12") '  -4.30%  '
Set-TextValue $ws.Range("D13") '19.84'
Set-TextValue $ws.Range("E13") '  -1.46%  '
Set-TextValue $ws.Range("D14") '7.62'
Set-TextValue $ws.Range("E14") '  -2.83%  '
Set-TextValue $ws.Range("D15") '3.220.44'
Set-TextValue $ws.Range("E15") '  -2.01%  '
Set-TextValue $ws.Range("D16") '2.759.15'
Set-TextValue $ws.Range("E16") '  -3.39%  '
Set-TextValue $ws.Range("E17") '  -0.78%  '
Set-TextValue $ws.Range("D18") '51.668.87'
Set-TextValue $ws.Range("E18") '  -1.14%  '
Set-TextValue $ws.Range("D19") '7.65'
Set-TextValue $ws.Range("E19") '  +0.78%  '
Set-TextValue $ws.Range("E20") '  -2.26%  '
Set-TextValue $ws.Range("D21") '13.25'
Set-TextValue $ws.Range("E21") '  -2.32%  '
Set-TextValue $ws.Range("E22") '  -3.44%  '
Set-TextValue $ws.Range("D23") '70.15'
Set-TextValue $ws.Range("E23") '  -0.82%  '
Set-TextValue $ws.Range("D24") '268.50'
Set-TextValue $ws.Range("E24") '  -1.73%  '
Set-TextValue $ws.Range("D25") '2.76'
Set-TextValue $ws.Range("E25") '  -3.05%  '
Set-TextValue $ws.Range("D26") '26.37'
Set-TextValue $ws.Range("E26") '  -3.35%  '
Set-TextValue $ws.Range("E27") '  +0.13%  '
Set-TextValue $ws.Range("E28") '  +15.27%  '
Set-TextValue $ws.Range("D29") '10.19'
Set-TextValue $ws.Range("E29") '  -1.93%  '
Set-TextValue $ws.Range("D30") '2.27'
Set-TextValue $ws.Range("E30") '  +0.61%  '
Set-TextValue $ws.Range("B31") 'InjectiveProtocol'
Set-TextValue $ws.Range("C31") 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D31") '35.43'
Set-TextValue $ws.Range("E31") '  -0.61%  '
Set-TextValue $ws.Range("B32") 'Filecoin'
Set-TextValue $ws.Range("C32") 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D32") '6.18'
Set-TextValue $ws.Range("E32") '  +4.27%  '
Set-TextValue $ws.Range("D33") '52.09'
Set-TextValue $ws.Range("E33") '  -0.49%  '
Set-TextValue $ws.Range("D34") '0.0444'
Set-TextValue $ws.Range("E34") '  -8.28%  '
Set-TextValue $ws.Range("D35") '0.0840'
Set-TextValue $ws.Range("E35") '  -2.03%  '
Set-TextValue $ws.Range("D36") '5.21'
Set-TextValue $ws.Range("E36") '  -8.07%  '
Set-TextValue $ws.Range("D37") '0.999'
Set-TextValue $ws.Range("E37") '  -0.10%  '
Set-TextValue $ws.Range("D38") '18.83'
Set-TextValue $ws.Range("E38") '  +1.26%  '
Set-TextValue $ws.Range("E39") '  -5.15%  '
Set-TextValue $ws.Range("E40") '  -4.88%  '
Set-TextValue $ws.Range("E41") '  -3.32%  '
Set-TextValue $ws.Range("E42") '  -1.57%  '
Set-TextValue $ws.Range("E43") '  -3.55%  '
Set-TextValue $ws.Range("D44") '119.72'
Set-TextValue $ws.Range("E44") '  -5.72%  '
Set-TextValue $ws.Range("D45") '21.90'
Set-TextValue $ws.Range("E45") '  -6.33%  '
Set-TextValue $ws.Range("D46") '2.095.91'
Set-TextValue $ws.Range("E46") '  -0.22%  '
Set-TextValue $ws.Range("D47") '3.26'
Set-TextValue $ws.Range("E47") '  -3.66%  '
Set-TextValue $ws.Range("D48") '2.32'
Set-TextValue $ws.Range("E48") '  +0.77%  '
Set-TextValue $ws.Range("E49") '  -4.50%  '
Set-TextValue $ws.Range("D50") '5.56'
Set-TextValue $ws.Range("E50") '  -6.33%  '
Set-TextValue $ws.Range("B51") 'Algorand'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D51") '0.191'
Set-TextValue $ws.Range("E51") '  -2.89%  '
